$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 819.4286
$ws.Range("I28").Value = 743.75
$ws.Range("J28").Value = 920.3333
$ws.Range("K28").Value = 743.75
$ws.Range("L28").Value = 920.3333
$ws.Range("M28").Value = -258.75
$ws.Range("N28").Value = -1890.3333

$ws.Range("H39").Value = 983.4286
$ws.Range("I39").Value = 21.2
$ws.Range("J39").Value = 3389
$ws.Range("K39").Value = 63.59999999999999
$ws.Range("L39").Value = 10167
$ws.Range("M39").Value = 232.4
$ws.Range("N39").Value = -10759

$ws.Range("H98").Value = 4259.647
$ws.Range("I98").Value = 3563.2
$ws.Range("K98").Value = 3563.2
$ws.Range("M98").Value = -2065.2

$ws.Range("H122").Value = 4259.647
$ws.Range("I122").Value = 3563.2
$ws.Range("K122").Value = 10689.6
$ws.Range("M122").Value = -8239.599999999999

$ws.Range("H137").Value = 4758.88
$ws.Range("I137").Value = 1571.5238
$ws.Range("J137").Value = 21492.5
$ws.Range("K137").Value = 4714.5714
$ws.Range("L137").Value = 64477.5
$ws.Range("M137").Value = -2164.5714
$ws.Range("N137").Value = -69577.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4091.6924
$ws.Range("I2").Value = 3020.3333
$ws.Range("J2").Value = 6502.25
$ws.Range("K2").Value = 3020.3333
$ws.Range("L2").Value = 6502.25
$ws.Range("M2").Value = -2907.3333
$ws.Range("N2").Value = -6728.25

$ws.Range("H116").Value = 4091.6924
$ws.Range("I116").Value = 3020.3333
$ws.Range("J116").Value = 6502.25
$ws.Range("K116").Value = 3020.3333
$ws.Range("L116").Value = 6502.25
$ws.Range("M116").Value = -726.3332999999998
$ws.Range("N116").Value = -11090.25

$ws.Range("H122").Value = 22225470
$ws.Range("I122").Value = 2556.3333
$ws.Range("K122").Value = 7668.999899999999
$ws.Range("M122").Value = -5218.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4091.6924
$ws.Range("I3").Value = 3020.3333
$ws.Range("J3").Value = 6502.25
$ws.Range("K3").Value = 3020.3333
$ws.Range("L3").Value = 6502.25
$ws.Range("M3").Value = -2906.3333
$ws.Range("N3").Value = -6730.25

$ws.Range("H107").Value = 2958.3845
$ws.Range("I107").Value = 2371.6667
$ws.Range("K107").Value = 2371.6667
$ws.Range("M107").Value = -451.6667000000002

$ws.Range("H134").Value = 2647.7058
$ws.Range("J134").Value = 1581
$ws.Range("L134").Value = 4743
$ws.Range("N134").Value = -9813

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 38890.25
$ws.Range("J48").Value = 38890.25
$ws.Range("L48").Value = 38890.25
$ws.Range("N48").Value = -39842.25

$ws.Range("H99").Value = 7453.7046
$ws.Range("I99").Value = 3936.7896
$ws.Range("J99").Value = 10126.56
$ws.Range("K99").Value = 3936.7896
$ws.Range("L99").Value = 10126.56
$ws.Range("M99").Value = -2438.7896
$ws.Range("N99").Value = -13122.56

$ws.Range("H102").Value = 58241
$ws.Range("J102").Value = 58241
$ws.Range("L102").Value = 58241
$ws.Range("N102").Value = -63109

$ws.Range("H105").Value = 11358.1
$ws.Range("I105").Value = 1303.3334
$ws.Range("K105").Value = 1303.3334
$ws.Range("M105").Value = 443.6666

$ws.Range("H126").Value = 7453.7046
$ws.Range("I126").Value = 3936.7896
$ws.Range("J126").Value = 10126.56
$ws.Range("K126").Value = 11810.3688
$ws.Range("L126").Value = 30379.68
$ws.Range("M126").Value = -9340.3688
$ws.Range("N126").Value = -35319.68

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 500
$ws.Range("I57").Value = 500
$ws.Range("K57").Value = 1500
$ws.Range("M57").Value = -941

$ws.Range("H121").Value = 1738.1666
$ws.Range("I121").Value = 1357.25
$ws.Range("K121").Value = 4071.75
$ws.Range("M121").Value = -2761.75

$ws.Range("H131").Value = 1055.3684
$ws.Range("I131").Value = 848.3333
$ws.Range("J131").Value = 1831.75
$ws.Range("K131").Value = 2544.9999
$ws.Range("L131").Value = 5495.25
$ws.Range("M131").Value = 2495.0001
$ws.Range("N131").Value = -15575.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 51126.25
$ws.Range("J26").Value = 51252.5
$ws.Range("L26").Value = 51252.5
$ws.Range("N26").Value = -51812.5

$ws.Range("H49").Value = 20185
$ws.Range("I49").Value = 18981.25
$ws.Range("K49").Value = 18981.25
$ws.Range("M49").Value = -18797.25

$ws.Range("H50").Value = 51126.25
$ws.Range("J50").Value = 51252.5
$ws.Range("L50").Value = 51252.5
$ws.Range("N50").Value = -52248.5

$ws.Range("H97").Value = 1288.9412
$ws.Range("I97").Value = 1354.1818
$ws.Range("K97").Value = 1354.1818
$ws.Range("M97").Value = -858.1818000000001

$ws.Range("H102").Value = 4634.0713
$ws.Range("I102").Value = 3151.5
$ws.Range("J102").Value = 5746
$ws.Range("K102").Value = 3151.5
$ws.Range("L102").Value = 5746
$ws.Range("M102").Value = -1529.5
$ws.Range("N102").Value = -8990

$ws.Range("H113").Value = 5313.2666
$ws.Range("I113").Value = 4600
$ws.Range("J113").Value = 5491.5835
$ws.Range("K113").Value = 4600
$ws.Range("L113").Value = 5491.5835
$ws.Range("M113").Value = -2430
$ws.Range("N113").Value = -9831.583500000001

$ws.Range("H122").Value = 17244010
$ws.Range("I122").Value = 2826.0625
$ws.Range("J122").Value = 38463930
$ws.Range("K122").Value = 8478.1875
$ws.Range("L122").Value = 115391790
$ws.Range("M122").Value = -6028.1875
$ws.Range("N122").Value = -115396690

$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2233.4119
$ws.Range("I16").Value = 1247.625
$ws.Range("J16").Value = 3109.6667
$ws.Range("K16").Value = 1247.625
$ws.Range("L16").Value = 3109.6667
$ws.Range("M16").Value = -1077.625
$ws.Range("N16").Value = -3449.6667

$ws.Range("H40").Value = 2885.75
$ws.Range("I40").Value = 3821.75
$ws.Range("J40").Value = 1949.75
$ws.Range("K40").Value = 3821.75
$ws.Range("L40").Value = 1949.75
$ws.Range("M40").Value = -3685.75
$ws.Range("N40").Value = -2221.75

$ws.Range("H46").Value = 1578.2759

$ws.Range("H50").Value = 34055.668
$ws.Range("J50").Value = 34055.668
$ws.Range("L50").Value = 34055.668
$ws.Range("N50").Value = -35329.668

$ws.Range("H55").Value = 870.625
$ws.Range("I55").Value = 728
$ws.Range("J55").Value = 981.55554
$ws.Range("K55").Value = 728
$ws.Range("L55").Value = 981.55554
$ws.Range("M55").Value = -555
$ws.Range("N55").Value = -1327.55554

$ws.Range("H56").Value = 44495
$ws.Range("J56").Value = 44495
$ws.Range("L56").Value = 44495
$ws.Range("N56").Value = -45877

$ws.Range("H82").Value = 3199.375
$ws.Range("J82").Value = 2560.6
$ws.Range("L82").Value = 2560.6
$ws.Range("N82").Value = -3282.6

$ws.Range("H85").Value = 3199.375
$ws.Range("J85").Value = 2560.6
$ws.Range("L85").Value = 2560.6
$ws.Range("N85").Value = -5056.6

$ws.Range("H93").Value = 524435.1
$ws.Range("I93").Value = 1965.3478
$ws.Range("K93").Value = 1965.3478
$ws.Range("M93").Value = -717.3478

$ws.Range("H122").Value = 3052790.5
$ws.Range("I122").Value = 3629.1292
$ws.Range("J122").Value = 12505190
$ws.Range("K122").Value = 10887.3876
$ws.Range("L122").Value = 37515570
$ws.Range("M122").Value = -8437.3876
$ws.Range("N122").Value = -37520470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 899.5
$ws.Range("I10").Value = 799
$ws.Range("J10").Value = 1000
$ws.Range("K10").Value = 799
$ws.Range("L10").Value = 1000
$ws.Range("M10").Value = -630
$ws.Range("N10").Value = -1338

$ws.Range("H13").Value = 1603.3334
$ws.Range("I13").Value = 1603.3334
$ws.Range("K13").Value = 1603.3334
$ws.Range("M13").Value = -1463.3334

$ws.Range("H14").Value = 2591.2144
$ws.Range("J14").Value = 2640.7693
$ws.Range("L14").Value = 2640.7693
$ws.Range("N14").Value = -2976.7693

$ws.Range("H17").Value = 484.66666
$ws.Range("I17").Value = 484.66666
$ws.Range("K17").Value = 484.66666
$ws.Range("M17").Value = -312.66666

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws.Range("H47").Value = 37495
$ws.Range("J47").Value = 37495
$ws.Range("L47").Value = 37495
$ws.Range("N47").Value = -38639

$ws.Range("H58").Value = 41322
$ws.Range("I58").Value = 39735.5
$ws.Range("J58").Value = 44495
$ws.Range("K58").Value = 39735.5
$ws.Range("L58").Value = 44495
$ws.Range("M58").Value = -39427.5
$ws.Range("N58").Value = -45111

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").ClearContents()

$ws.Range("H70").Value = 32998
$ws.Range("J70").Value = 32998
$ws.Range("L70").Value = 32998
$ws.Range("N70").Value = -33628

$ws.Range("H73").Value = 32998
$ws.Range("J73").Value = 32998
$ws.Range("L73").Value = 32998
$ws.Range("N73").Value = -35182

$ws.Range("H98").Value = 7777
$ws.Range("J98").Value = 7777
$ws.Range("L98").Value = 7777
$ws.Range("N98").Value = -13767
